# Update "provincias_spain" city stats table (sheet "Ciudades") with the
# latest data refresh: numbers for most provinces move, and because the
# table stays sorted by "Casos totales" descending, several rows swap the
# city name shown in column A. The "last updated" banner in A1 also shifts
# from 14:05 to 14:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Map of row -> @{ column letter = new value } for every cell that changed.
$updates = @{
    1 = @{ "A" = 'Datos actualizados a 15 de Mayo de 2020 a las 14:35' }
    12 = @{ "A" = 'Ciudad Real'; "B" = 6428; "C" = 1996; "D" = 3377; "E" = 1055 }
    13 = @{ "A" = 'Valencia/Valencia'; "B" = 5565; "C" = 4828; "D" = 2767; "E" = 687 }
    15 = @{ "B" = 5170; "C" = 2676; "D" = 1840; "E" = 654 }
    18 = @{ "A" = 'Valladolid'; "B" = 4348; "C" = 1554; "D" = 2436; "E" = 358 }
    19 = @{ "A" = 'Salamanca'; "B" = 4136; "C" = 1173; "D" = 2608; "E" = 355 }
    20 = @{ "A" = 'La Rioja'; "B" = 4016; "C" = 2927; "D" = 741; "E" = 348 }
    21 = @{ "A" = 'Malaga'; "B" = 3983; "C" = 2161; "D" = 1548; "E" = 274 }
    22 = @{ "A" = 'Toledo'; "B" = 3846; "C" = 1929; "D" = 1165; "E" = 752 }
    23 = @{ "A" = 'Alacant/Alicante'; "B" = 3768; "C" = 3541; "D" = 1938; "E" = 480 }
    24 = @{ "B" = 3764; "D" = 1893; "E" = 508 }
    25 = @{ "B" = 3551; "C" = 1594; "D" = 1553; "E" = 404 }
    26 = @{ "B" = 3404; "C" = 861; "D" = 2343 }
    28 = @{ "B" = 3091; "C" = 1742; "D" = 1075; "E" = 274 }
    29 = @{ "B" = 3063; "C" = 2480; "D" = 305; "E" = 278 }
    31 = @{ "B" = 2728; "C" = 897; "D" = 1626 }
    34 = @{ "B" = 2278; "C" = 393; "D" = 1766 }
    36 = @{ "A" = 'Caceres'; "B" = 1973; "C" = 1505; "D" = 66; "E" = 402 }
    37 = @{ "A" = 'A Coruña'; "B" = 1969; "C" = 333; "D" = 1788; "E" = 67 }
    38 = @{ "A" = 'Avila'; "B" = 1917; "C" = 618; "D" = 1166; "E" = 133 }
    39 = @{ "A" = 'Jaen'; "B" = 1751; "C" = 1121; "D" = 457; "E" = 173 }
    40 = @{ "A" = 'Cordoba'; "B" = 1682; "C" = 1331; "D" = 246; "E" = 105 }
    41 = @{ "A" = 'Pontevedra'; "B" = 1536; "C" = 333; "D" = 1411; "E" = 30 }
    42 = @{ "A" = 'Tenerife'; "B" = 1532; "C" = 897; "D" = 529; "E" = 106 }
    43 = @{ "A" = 'Murcia'; "B" = 1508; "C" = 1782; "D" = 0; "E" = 139 }
    44 = @{ "A" = 'Cadiz'; "B" = 1480; "C" = 517; "D" = 821; "E" = 142 }
    45 = @{ "A" = 'Castello/Castellon'; "B" = 1475; "C" = 1339; "D" = 699; "E" = 206 }
    46 = @{ "D" = 640; "E" = 250 }
    47 = @{ "B" = 1216; "D" = 320; "E" = 305 }
    48 = @{ "B" = 1191; "C" = 326; "D" = 784; "E" = 81 }
    49 = @{ "B" = 1090; "C" = 441; "D" = 550 }
    50 = @{ "B" = 962; "C" = 1082 }
    51 = @{ "B" = 957; "C" = 319; "D" = 553; "E" = 85 }
    53 = @{ "B" = 695; "C" = 452; "D" = 193 }
    54 = @{ "B" = 646; "C" = 365; "D" = 198 }
    56 = @{ "B" = 520; "C" = 344; "D" = 128 }
}

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $ws.Range("$col$row").Value = $rowUpdates[$col]
    }
}

